# Apply the "MaxPowerLevel, MaxStageLevel int keys + DamageRateTable" edit.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) GlobalConstantIntTable (1st sheet): drop GlobalConstant_MaxRank,
#    add MaxPowerLevel / MaxStageLevel (both int 10).
# ---------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item(1)

$wsInt.Range("A1").Value = "id|String"
$wsInt.Range("B1").Value = "value|Int"

$wsInt.Range("A2").Value = "MaxPowerLevel"
$wsInt.Range("B2").Value = 10

$wsInt.Range("A3").Value = "MaxStageLevel"
$wsInt.Range("B3").Value = 10

# ---------------------------------------------------------------------
# 2) GlobalConstantFloatTable (2nd sheet): unchanged content
#    (id|String / value|Float / SpDecreaseRate / 0.8).
# ---------------------------------------------------------------------
$wsFloat = $wb.Worksheets.Item(2)

$wsFloat.Range("A1").Value = "id|String"
$wsFloat.Range("B1").Value = "value|Float"

$wsFloat.Range("A2").Value = "SpDecreaseRate"
$wsFloat.Range("B2").Value = 0.8

# ---------------------------------------------------------------------
# 3) New sheet "DamageRateTable" appended at the end.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsDmg = $wb.Worksheets.Add($null, $lastSheet)
$wsDmg.Name = "DamageRateTable"

$wsDmg.Range("A1").Value = "id|String"
$wsDmg.Range("B1").Value = "number|Int"
$wsDmg.Range("C1").Value = "rate|Float!"

$wsDmg.Range("A2").Value = "Ricochet"
$wsDmg.Range("B2").Value = 1
$wsDmg.Range("C2").Value = "1,0.7"

$wsDmg.Range("A3").Value = "Ricochet"
$wsDmg.Range("B3").Value = 2
$wsDmg.Range("C3").Value = "1,0.7,0.49"

$wsDmg.Range("A4").Value = "Ricochet"
$wsDmg.Range("B4").Value = 3
$wsDmg.Range("C4").Value = "1,0.7,0.2"

$wsDmg.Range("A5").Value = "BounceWallQuad"
$wsDmg.Range("B5").Value = 1
$wsDmg.Range("C5").Value = "1,0.5"

$wsDmg.Range("A6").Value = "BounceWallQuad"
$wsDmg.Range("B6").Value = 2
$wsDmg.Range("C6").Value = "1,0.5,0.25"

$wsDmg.Range("A7").Value = "MonsterThrough"
$wsDmg.Range("B7").Value = 1
$wsDmg.Range("C7").Value = "1,0.66"

$wsDmg.Range("A8").Value = "MonsterThrough"
$wsDmg.Range("B8").Value = 1
$wsDmg.Range("C8").Value = "1,0.66,0.3"

$wsDmg.Range("A9").Value = "Repeat"
$wsDmg.Range("B9").Value = 1
$wsDmg.Range("C9").Value = "1,0.9"

$wsDmg.Range("A10").Value = "Repeat"
$wsDmg.Range("B10").Value = 2
$wsDmg.Range("C10").Value = "1,0.9,0.81"

$wsDmg.Range("A11").Value = "Repeat"
$wsDmg.Range("B11").Value = 3
$wsDmg.Range("C11").Value = "1.0.9,0.72"

$wsDmg.Range("A12").Value = "Parallel"
$wsDmg.Range("B12").Value = 1
$wsDmg.Range("C12").Value = 0.75

$wsDmg.Range("A13").Value = "Parallel"
$wsDmg.Range("B13").Value = 2
$wsDmg.Range("C13").Value = 0.5625

$wsDmg.Range("A14").Value = "CircularSector"
$wsDmg.Range("B14").Value = 1
$wsDmg.Range("C14").Value = 0.8

$wsDmg.Range("A15").Value = "WallThrough"
$wsDmg.Range("C15").Value = 0.75

# ---------------------------------------------------------------------
# 4) Restore the active tab to the first sheet (GlobalConstantIntTable),
#    matching the saved workbookView (activeTab back to default / 0).
# ---------------------------------------------------------------------
$wsInt.Activate()
